$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra prompt-variant columns D through L (keep only A, B, C)
$ws.Range("D1:L1").ClearContents()
$ws.Columns("D:L").Delete()

# Update C1 with the new composition-check prompt
$ws.Range("C1").Value = "In a textual domain description, there is one statement (statement 1) that says: Statement 1: {statement1}.  We also have a second statement (statement 2) that says:  Statement 2: {statement2}.  Please, analyze statement 2 to determine whether {source} is composed of {target}.  Begin your response by providing reasoning, and conclude:  - 'Conclusion:Yes' if {source} is composed of {target}  - 'Conclusion:No' if {source} is not composed of {target}, - 'Conclusion:Not Sure' if statement 2 does not provide enough details to determine whether  {source} is composed of {target}."

# Row 2
$ws.Range("A2").Value = "The game is an intergalactic murder mystery, which emulates a board game scenario in which players spin a wheel to determine a randomly selected number of moves by which they travel through space to various planets."
$ws.Range("B2").Value = "Each board is made up of cells"
$ws.Range("C2").Value = "Statement 2 describes the game as an intergalactic murder mystery that mirrors a board game. However, it does not explicitly mention what the board is composed of. It only tells us about the gameplay involving players spinning a wheel to decide their moves. Nothing is stated about whether the board consists of cells or not.`nConclusion: Not Sure"

# Row 3
$ws.Range("A3").Value = "These entry points are graphically depicted on the game as a diagonal line off a square on the board leading into a planet's atmosphere."
$ws.Range("B3").Value = "Each board is made up of cells"
$ws.Range("C3").Value = "Statement 2 discusses the graphical depiction of entry points on the game board, including the feature of a diagonal line which leads off a square on the board, leading into a planet's atmosphere. However, it does not specifically mention or describe the components of the board itself, for instance - whether the board is composed of cells as suggested in statement 1. `nConclusion: Not Sure"

# Row 4
$ws.Range("A4").Value = "A player visits a planet by making legal moves on the board until a cell adjacent to an atmosphere entry point is reached."
$ws.Range("B4").Value = "Each board is made up of cells"
$ws.Range("C4").Value = "Statement 2 indirectly suggests that the board is composed of cells by implying players make moves onto cells of the board. Although it does not explicitly state this, it does harmonize with statement 1 that says each board is composed of cells. Therefore, putting context provided in statement 1 and implied content of statement 2 together, we infer that the board is composed of cells.`nConclusion: Yes"

# Row 5
$ws.Range("A5").Value = "If an announced hypothesis is incorrect, the player loses the game and cannot pose hypotheses any longer or make moves on the board but must continue refuting the hypotheses of other players."
$ws.Range("B5").Value = "Each board is made up of cells"
$ws.Range("C5").Value = "Statement 2 discusses the rules of a hypothetical game, specifically the consequences a player faces if they make an incorrect hypothesis. However, it does not provide any information about the composition or structure of the board that is used in the game. Therefore, we cannot deduce from statement 2 whether or not the board is made up of cells.`nConclusion: Not Sure"

# Row 6
$ws.Range("A6").Value = "Each planet has between one and three predefined atmospheric entry points."
$ws.Range("B6").Value = "Each planet is made up of one to three planets"
$ws.Range("C6").Value = "Statement 2 specifies that each planet has between one and three predefined atmospheric entry points. However, the phrasing does not suggest that the planet is composed of these entry points. Rather, it implies that these entry points are simply located on or around the planet. Their existence does not necessarily mean that they form a part of the planet's essential composition.`nConclusion: No."

# Values containing embedded line breaks can cause the engine to auto-expand
# row height; restore rows to their default (non-custom) height.
$ws.Rows("2:6").AutoFit()
